$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf7"
$ws.Range("C2").Value = "Fgfr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.79659733333333
$ws.Range("H2").Value = 41.389792
$ws.Range("I2").Value = 0.9485830781324925
$ws.Range("J2").Value = 0.9485830781324925
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.4871643333333333
$ws.Range("N2").Value = 1.461493
$ws.Range("O2").Value = 0.1914458468300136
$ws.Range("P2").Value = 0.1914458468300136
$ws.Range("Q2").Value = 6.721210142161778
$ws.Range("R2").Value = 60.49089127945599
$ws.Range("S2").Value = 0.181602290681696
$ws.Range("T2").Value = 0.181602290681696

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf7"
$ws.Range("C3").Value = "Fgfr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.79659733333333
$ws.Range("H3").Value = 41.389792
$ws.Range("I3").Value = 0.9485830781324925
$ws.Range("J3").Value = 0.9485830781324925
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.864751
$ws.Range("N3").Value = 5.594253
$ws.Range("O3").Value = 0.7328098752209857
$ws.Range("P3").Value = 0.7328098752209857
$ws.Range("Q3").Value = 25.72721867393066
$ws.Range("R3").Value = 231.544968065376
$ws.Range("S3").Value = 0.6951310471230104
$ws.Range("T3").Value = 0.6951310471230104

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf7"
$ws.Range("C4").Value = "Fgfr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.79659733333333
$ws.Range("H4").Value = 41.389792
$ws.Range("I4").Value = 0.9485830781324925
$ws.Range("J4").Value = 0.9485830781324925
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1927433333333333
$ws.Range("N4").Value = 0.57823
$ws.Range("O4").Value = 0.07574427794900063
$ws.Range("P4").Value = 0.07574427794900063
$ws.Range("Q4").Value = 2.659202158684445
$ws.Range("R4").Value = 23.93281942816
$ws.Range("S4").Value = 0.0718497403277861
$ws.Range("T4").Value = 0.0718497403277861

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf7"
$ws.Range("C5").Value = "Fgfr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7478296666666666
$ws.Range("H5").Value = 2.243489
$ws.Range("I5").Value = 0.05141692186750751
$ws.Range("J5").Value = 0.05141692186750751
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.4871643333333333
$ws.Range("N5").Value = 1.461493
$ws.Range("O5").Value = 0.1914458468300136
$ws.Range("P5").Value = 0.1914458468300136
$ws.Range("Q5").Value = 0.3643159410085555
$ws.Range("R5").Value = 3.278843469077
$ws.Range("S5").Value = 0.00984355614831762
$ws.Range("T5").Value = 0.00984355614831762

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf7"
$ws.Range("C6").Value = "Fgfr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7478296666666666
$ws.Range("H6").Value = 2.243489
$ws.Range("I6").Value = 0.05141692186750751
$ws.Range("J6").Value = 0.05141692186750751
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.864751
$ws.Range("N6").Value = 5.594253
$ws.Range("O6").Value = 0.7328098752209857
$ws.Range("P6").Value = 0.7328098752209857
$ws.Range("Q6").Value = 1.394516118746333
$ws.Range("R6").Value = 12.550645068717
$ws.Range("S6").Value = 0.03767882809797535
$ws.Range("T6").Value = 0.03767882809797535

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf7"
$ws.Range("C7").Value = "Fgfr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7478296666666666
$ws.Range("H7").Value = 2.243489
$ws.Range("I7").Value = 0.05141692186750751
$ws.Range("J7").Value = 0.05141692186750751
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1927433333333333
$ws.Range("N7").Value = 0.57823
$ws.Range("O7").Value = 0.07574427794900063
$ws.Range("P7").Value = 0.07574427794900063
$ws.Range("Q7").Value = 0.1441391827188889
$ws.Range("R7").Value = 1.29725264447
$ws.Range("S7").Value = 0.003894537621214538
$ws.Range("T7").Value = 0.003894537621214538

